$d = $word.ActiveDocument

$replacements = @(
    @("810÷9=", "674÷7="),
    @("575÷8=", "218÷4="),
    @("410÷9=", "777÷9="),
    @("158÷2=", "777÷4="),
    @("514÷7=", "748÷6="),
    @("378÷7=", "847÷4="),
    @("410÷8=", "428÷7="),
    @("274÷8=", "646÷2="),
    @("777÷6=", "188÷2="),
    @("156÷5=", "154÷6="),
    @("650÷6=", "242÷4="),
    @("997÷6=", "314÷6="),
    @("607÷5=", "306÷4="),
    @("319÷3=", "855÷6="),
    @("531÷9=", "329÷7="),
    @("377÷5=", "583÷8="),
    @("585÷9=", "849÷6="),
    @("313÷8=", "741÷9="),
    @("650÷7=", "662÷4="),
    @("830÷4=", "196÷6="),
    @("136÷4=", "880÷9="),
    @("246÷6=", "860÷2="),
    @("965÷9=", "666÷7="),
    @("923÷7=", "626÷7="),
    @("536÷6=", "364÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
